# A new weekly price record for "Papa" / "Asterix" / "1a (cosecha)" was
# published at the top of the data block (row 8). All the existing records
# that used to occupy rows 8-68 shift down one row (to rows 9-69); the
# worksheet's used range grows from A1:R68 to A1:R69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (rows 8 and below) down by one row.
$ws.Rows.Item(8).Insert()

# Populate the newly freed row 8 with the new market entry.
$ws.Cells.Item(8, 1).Value  = 1
$ws.Cells.Item(8, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value  = 44670
$ws.Cells.Item(8, 5).Value  = 15
$ws.Cells.Item(8, 6).Value  = 100114001
$ws.Cells.Item(8, 7).Value  = "Papa"
$ws.Cells.Item(8, 8).Value  = "Asterix"
$ws.Cells.Item(8, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 8500
$ws.Cells.Item(8, 12).Value = 9000
$ws.Cells.Item(8, 13).Value = 8750
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(8, 16).Value = 350
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
